# Auto-generated edit script applying the Siren_Profits.xlsx diff
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H, I, J, K, L, M, N) for specific rows across 8 worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 15101.25
$ws.Range("I28").Value = 18466.334
$ws.Range("K28").Value = 18466.334
$ws.Range("M28").Value = -17981.334

$ws.Range("H74").Value = 4757.75
$ws.Range("I74").Value = 4367
$ws.Range("K74").Value = 4367
$ws.Range("M74").Value = -3431

$ws.Range("H77").Value = 4757.75
$ws.Range("I77").Value = 4367
$ws.Range("K77").Value = 21835
$ws.Range("M77").Value = -17155

$ws.Range("H111").Value = 4538.4
$ws.Range("I111").Value = 4538.4
$ws.Range("K111").Value = 13615.2
$ws.Range("M111").Value = -10548.2

$ws.Range("H113").Value = 11806.529
$ws.Range("I113").Value = 18764.777
$ws.Range("K113").Value = 18764.777
$ws.Range("M113").Value = -15510.777

$ws.Range("H116").Value = 834772.25
$ws.Range("I116").Value = 1414875
$ws.Range("K116").Value = 1414875
$ws.Range("M116").Value = -1411433

$ws.Range("H121").Value = 2000
$ws.Range("J121").Value = 2000
$ws.Range("L121").Value = 6000
$ws.Range("N121").Value = -9494

$ws.Range("H135").Value = 11763.637
$ws.Range("I135").Value = 12711.111
$ws.Range("J135").Value = 7500
$ws.Range("K135").Value = 114399.999
$ws.Range("L135").Value = 67500
$ws.Range("M135").Value = -111864.999
$ws.Range("N135").Value = -72570

$ws.Range("H141").Value = 5105.8335
$ws.Range("I141").Value = 4238.9473
$ws.Range("K141").Value = 12716.8419
$ws.Range("M141").Value = -7536.841899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4034.4333
$ws.Range("I32").Value = 4176.25
$ws.Range("K32").Value = 4176.25
$ws.Range("M32").Value = -3889.25

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 20096.5
$ws.Range("J95").Value = 20096.5
$ws.Range("L95").Value = 20096.5
$ws.Range("N95").Value = -25588.5

$ws.Range("H105").Value = 2023.375
$ws.Range("I105").Value = 1687.1818
$ws.Range("K105").Value = 1687.1818
$ws.Range("M105").Value = 59.81819999999993

$ws.Range("H107").Value = 3299.2917
$ws.Range("I107").Value = 3010.2222
$ws.Range("J107").Value = 4166.5
$ws.Range("K107").Value = 3010.2222
$ws.Range("L107").Value = 4166.5
$ws.Range("M107").Value = -1090.2222
$ws.Range("N107").Value = -8006.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3206.6743
$ws.Range("I31").Value = 1943.2609
$ws.Range("J31").Value = 4659.6
$ws.Range("K31").Value = 1943.2609
$ws.Range("L31").Value = 4659.6
$ws.Range("M31").Value = -1648.2609
$ws.Range("N31").Value = -5249.6

$ws.Range("H34").Value = 3206.6743
$ws.Range("I34").Value = 1943.2609
$ws.Range("J34").Value = 4659.6
$ws.Range("K34").Value = 1943.2609
$ws.Range("L34").Value = 4659.6
$ws.Range("M34").Value = -1741.2609
$ws.Range("N34").Value = -5063.6

$ws.Range("H86").Value = 9418.478999999999
$ws.Range("I86").Value = 9152.134
$ws.Range("J86").Value = 9917.875
$ws.Range("K86").Value = 9152.134
$ws.Range("L86").Value = 9917.875
$ws.Range("M86").Value = -8029.134
$ws.Range("N86").Value = -12163.875

$ws.Range("H89").Value = 9418.478999999999
$ws.Range("I89").Value = 9152.134
$ws.Range("J89").Value = 9917.875
$ws.Range("K89").Value = 45760.67
$ws.Range("L89").Value = 49589.375
$ws.Range("M89").Value = -40144.67
$ws.Range("N89").Value = -60821.375

$ws.Range("H99").Value = 1025844.8
$ws.Range("I99").Value = 1025844.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1025844.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1024346.8
$ws.Range("N99").Value = ""

$ws.Range("H107").Value = 55568852
$ws.Range("I107").Value = 71444936
$ws.Range("K107").Value = 71444936
$ws.Range("M107").Value = -71443016

$ws.Range("H122").Value = 767
$ws.Range("I122").Value = 756.1667
$ws.Range("J122").Value = 799.5
$ws.Range("K122").Value = 2268.5001
$ws.Range("L122").Value = 2398.5
$ws.Range("M122").Value = 181.4998999999998
$ws.Range("N122").Value = -7298.5

$ws.Range("H126").Value = 1025844.8
$ws.Range("I126").Value = 1025844.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3077534.4
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3075064.4
$ws.Range("N126").Value = ""

$ws.Range("H132").Value = 17862.559
$ws.Range("I132").Value = 6763.7036
$ws.Range("J132").Value = 60672.43
$ws.Range("K132").Value = 20291.1108
$ws.Range("L132").Value = 182017.29
$ws.Range("M132").Value = -17761.1108
$ws.Range("N132").Value = -187077.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""

$ws.Range("H94").Value = 89834100
$ws.Range("I94").Value = 1000000
$ws.Range("J94").Value = 112042620
$ws.Range("K94").Value = 1000000
$ws.Range("L94").Value = 112042620
$ws.Range("M94").Value = -999324
$ws.Range("N94").Value = -112043972

$ws.Range("H102").Value = 6421.7085
$ws.Range("I102").Value = 7374.7896
$ws.Range("J102").Value = 2800
$ws.Range("K102").Value = 7374.7896
$ws.Range("L102").Value = 2800
$ws.Range("M102").Value = -5752.7896
$ws.Range("N102").Value = -6044

$ws.Range("H107").Value = 571.9
$ws.Range("I107").Value = 502.9375
$ws.Range("K107").Value = 502.9375
$ws.Range("M107").Value = 1417.0625

$ws.Range("H122").Value = 18661
$ws.Range("I122").Value = 16646.53
$ws.Range("J122").Value = 27222.5
$ws.Range("K122").Value = 49939.59
$ws.Range("L122").Value = 81667.5
$ws.Range("M122").Value = -47489.59
$ws.Range("N122").Value = -86567.5

$ws.Range("H123").Value = 14649.95
$ws.Range("J123").Value = 14649.95
$ws.Range("L123").Value = 14649.95
$ws.Range("N123").Value = -19549.95

$ws.Range("H125").Value = 84775
$ws.Range("J125").Value = 84775
$ws.Range("L125").Value = 84775
$ws.Range("N125").Value = -89695

$ws.Range("H132").Value = 3116.8928
$ws.Range("I132").Value = 2229.7827
$ws.Range("K132").Value = 6689.348100000001
$ws.Range("M132").Value = -4159.348100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4697.143
$ws.Range("I16").Value = 3813.3333
$ws.Range("K16").Value = 3813.3333
$ws.Range("M16").Value = -3643.3333

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""

$ws.Range("H30").Value = 1533
$ws.Range("I30").Value = 377.33334
$ws.Range("J30").Value = 5000
$ws.Range("K30").Value = 377.33334
$ws.Range("L30").Value = 5000
$ws.Range("M30").Value = -269.33334
$ws.Range("N30").Value = -5216

$ws.Range("H55").Value = 862.65625
$ws.Range("I55").Value = 772.52
$ws.Range("K55").Value = 772.52
$ws.Range("M55").Value = -599.52

$ws.Range("H132").Value = 1530245.4
$ws.Range("I132").Value = 2293504.8
$ws.Range("J132").Value = 3726.5833
$ws.Range("K132").Value = 6880514.399999999
$ws.Range("L132").Value = 11179.7499
$ws.Range("M132").Value = -6877984.399999999
$ws.Range("N132").Value = -16239.7499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 10000
$ws.Range("K21").Value = 10000
$ws.Range("M21").Value = -9765

$ws.Range("H35").Value = 10000
$ws.Range("I35").Value = 10000
$ws.Range("K35").Value = 10000
$ws.Range("M35").Value = -9710

$ws.Range("H49").Value = 15500
$ws.Range("J49").Value = 15500
$ws.Range("L49").Value = 15500
$ws.Range("N49").Value = -15960

$ws.Range("H54").Value = 50000
$ws.Range("J54").Value = 50000
$ws.Range("L54").Value = 50000
$ws.Range("N54").Value = -51040

$ws.Range("H122").Value = 22604.393
$ws.Range("I122").Value = 2632.8635
$ws.Range("K122").Value = 7898.5905
$ws.Range("M122").Value = -5448.5905
